# Update Sheets via scheduled runner — price/profit refresh across worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4341.5454
$ws.Range("J17").Value = 4341.5454
$ws.Range("L17").Value = 13024.6362
$ws.Range("N17").Value = -13360.6362
$ws.Range("H61").Value = 299.2
$ws.Range("I61").Value = 299.2
$ws.Range("K61").Value = 897.5999999999999
$ws.Range("M61").Value = -725.5999999999999
$ws.Range("H130").Value = 43501.6
$ws.Range("J130").Value = 43501.6
$ws.Range("L130").Value = 43501.6
$ws.Range("N130").Value = -53541.6
$ws.Range("H132").Value = 23704.738
$ws.Range("I132").Value = 3310.1765
$ws.Range("J132").Value = 110381.625
$ws.Range("K132").Value = 9930.529500000001
$ws.Range("L132").Value = 331144.875
$ws.Range("M132").Value = -7400.529500000001
$ws.Range("N132").Value = -336204.875
$ws.Range("H137").Value = 3286.4844
$ws.Range("I137").Value = 930.125
$ws.Range("J137").Value = 4071.9375
$ws.Range("K137").Value = 2790.375
$ws.Range("L137").Value = 12215.8125
$ws.Range("M137").Value = -240.375
$ws.Range("N137").Value = -17315.8125
$ws.Range("H141").Value = 4539.3477
$ws.Range("I141").Value = 2615
$ws.Range("K141").Value = 7845
$ws.Range("M141").Value = -2665
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 862.55
$ws.Range("I2").Value = 793.4286
$ws.Range("J2").Value = 1346.4
$ws.Range("K2").Value = 793.4286
$ws.Range("L2").Value = 1346.4
$ws.Range("M2").Value = -680.4286
$ws.Range("N2").Value = -1572.4
$ws.Range("H32").Value = 35804.16
$ws.Range("I32").Value = 36137.668
$ws.Range("J32").Value = 32802.6
$ws.Range("K32").Value = 36137.668
$ws.Range("L32").Value = 32802.6
$ws.Range("M32").Value = -35850.668
$ws.Range("N32").Value = -33376.6
$ws.Range("H76").Value = 23644
$ws.Range("J76").Value = 23644
$ws.Range("L76").Value = 23644
$ws.Range("N76").Value = -24320
$ws.Range("H79").Value = 23644
$ws.Range("J79").Value = 23644
$ws.Range("L79").Value = 23644
$ws.Range("N79").Value = -25984
$ws.Range("H116").Value = 862.55
$ws.Range("I116").Value = 793.4286
$ws.Range("J116").Value = 1346.4
$ws.Range("K116").Value = 793.4286
$ws.Range("L116").Value = 1346.4
$ws.Range("M116").Value = 1500.5714
$ws.Range("N116").Value = -5934.4
$ws.Range("H134").Value = 52271.43
$ws.Range("J134").Value = 52271.43
$ws.Range("L134").Value = 52271.43
$ws.Range("N134").Value = -62411.43
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 862.55
$ws.Range("I3").Value = 793.4286
$ws.Range("J3").Value = 1346.4
$ws.Range("K3").Value = 793.4286
$ws.Range("L3").Value = 1346.4
$ws.Range("M3").Value = -679.4286
$ws.Range("N3").Value = -1574.4
$ws.Range("H126").Value = 50768
$ws.Range("J126").Value = 50768
$ws.Range("L126").Value = 50768
$ws.Range("N126").Value = -60648
$ws.Range("H137").Value = 53666.332
$ws.Range("J137").Value = 53666.332
$ws.Range("L137").Value = 53666.332
$ws.Range("N137").Value = -63866.332
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 238367.19
$ws.Range("I31").Value = 2258.1052
$ws.Range("J31").Value = 347783.6
$ws.Range("K31").Value = 2258.1052
$ws.Range("L31").Value = 347783.6
$ws.Range("M31").Value = -1963.1052
$ws.Range("N31").Value = -348373.6
$ws.Range("H34").Value = 238367.19
$ws.Range("I34").Value = 2258.1052
$ws.Range("J34").Value = 347783.6
$ws.Range("K34").Value = 2258.1052
$ws.Range("L34").Value = 347783.6
$ws.Range("M34").Value = -2056.1052
$ws.Range("N34").Value = -348187.6
$ws.Range("H124").Value = 26211.334
$ws.Range("J124").Value = 26211.334
$ws.Range("L124").Value = 26211.334
$ws.Range("N124").Value = -31121.334
$ws.Range("H125").Value = 49326
$ws.Range("J125").Value = 49326
$ws.Range("L125").Value = 49326
$ws.Range("N125").Value = -54246
$ws.Range("H131").Value = 38318
$ws.Range("J131").Value = 38318
$ws.Range("L131").Value = 38318
$ws.Range("N131").Value = -48398
$ws.Range("H132").Value = 66701.13
$ws.Range("I132").Value = 2639.7693
$ws.Range("J132").Value = 159234.22
$ws.Range("K132").Value = 7919.3079
$ws.Range("L132").Value = 477702.66
$ws.Range("M132").Value = -5389.3079
$ws.Range("N132").Value = -482762.66
$ws.Range("H134").Value = 425267.06
$ws.Range("I134").Value = 812.34784
$ws.Range("J134").Value = 1401512.9
$ws.Range("K134").Value = 2437.04352
$ws.Range("L134").Value = 4204538.699999999
$ws.Range("M134").Value = 97.95647999999983
$ws.Range("N134").Value = -4209608.699999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4637.04
$ws.Range("I113").Value = 8293.23
$ws.Range("J113").Value = 676.1667
$ws.Range("K113").Value = 24879.69
$ws.Range("L113").Value = 2028.5001
$ws.Range("M113").Value = -22709.69
$ws.Range("N113").Value = -6368.5001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 28868
$ws.Range("J118").Value = 28868
$ws.Range("L118").Value = 28868
$ws.Range("N118").Value = -32182
$ws.Range("H120").Value = 25283.666
$ws.Range("J120").Value = 25283.666
$ws.Range("L120").Value = 25283.666
$ws.Range("N120").Value = -34959.666
$ws.Range("H125").Value = 40996
$ws.Range("J125").Value = 40996
$ws.Range("L125").Value = 40996
$ws.Range("N125").Value = -45916
$ws.Range("H127").Value = 41996
$ws.Range("J127").Value = 41996
$ws.Range("L127").Value = 41996
$ws.Range("N127").Value = -51916
$ws.Range("H131").Value = 38986
$ws.Range("J131").Value = 38986
$ws.Range("L131").Value = 38986
$ws.Range("N131").Value = -49066
$ws.Range("H135").Value = 44239.8
$ws.Range("J135").Value = 44239.8
$ws.Range("L135").Value = 44239.8
$ws.Range("N135").Value = -54379.8
$ws.Range("H139").Value = 25567.111
$ws.Range("J139").Value = 25567.111
$ws.Range("L139").Value = 25567.111
$ws.Range("N139").Value = -35847.111
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2073.7334
$ws.Range("I93").Value = 2634
$ws.Range("J93").Value = 1933.6666
$ws.Range("K93").Value = 2634
$ws.Range("L93").Value = 1933.6666
$ws.Range("M93").Value = -1386
$ws.Range("N93").Value = -4429.6666
$ws.Range("H94").Value = 40282.9
$ws.Range("J94").Value = 40282.9
$ws.Range("L94").Value = 40282.9
$ws.Range("N94").Value = -41634.9
$ws.Range("H95").Value = 33598.2
$ws.Range("J95").Value = 33598.2
$ws.Range("L95").Value = 33598.2
$ws.Range("N95").Value = -39090.2
$ws.Range("H96").Value = 30664
$ws.Range("J96").Value = 30664
$ws.Range("L96").Value = 30664
$ws.Range("N96").Value = -36156
$ws.Range("H114").Value = 38382
$ws.Range("J114").Value = 38382
$ws.Range("L114").Value = 38382
$ws.Range("N114").Value = -47060
$ws.Range("H117").Value = 45388
$ws.Range("J117").Value = 45388
$ws.Range("L117").Value = 45388
$ws.Range("N117").Value = -54566
$ws.Range("H123").Value = 32872.668
$ws.Range("J123").Value = 32872.668
$ws.Range("L123").Value = 32872.668
$ws.Range("N123").Value = -42672.668
$ws.Range("H131").Value = 43318
$ws.Range("J131").Value = 43318
$ws.Range("L131").Value = 43318
$ws.Range("N131").Value = -53398
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 28906.5
$ws.Range("J27").Value = 28906.5
$ws.Range("L27").Value = 28906.5
$ws.Range("N27").Value = -29044.5
$ws.Range("H97").Value = 38134.25
$ws.Range("J97").Value = 38134.25
$ws.Range("L97").Value = 38134.25
$ws.Range("N97").Value = -40116.25
$ws.Range("H127").Value = 31561.334
$ws.Range("J127").Value = 31561.334
$ws.Range("L127").Value = 31561.334
$ws.Range("N127").Value = -41481.334
$ws.Range("H129").Value = 28694.666
$ws.Range("J129").Value = 28694.666
$ws.Range("L129").Value = 28694.666
$ws.Range("N129").Value = -38694.666
$ws.Range("H133").Value = 102728.5
$ws.Range("J133").Value = 102728.5
$ws.Range("L133").Value = 102728.5
$ws.Range("N133").Value = -112848.5
$ws.Range("H137").Value = 56359.75
$ws.Range("J137").Value = 56359.75
$ws.Range("L137").Value = 56359.75
$ws.Range("N137").Value = -66559.75
